$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 303, shifting existing rows 303-341 down to 304-342.
$ws.Rows.Item(303).Insert()

# Populate the newly inserted row 303 with the new weekly record.
$ws.Cells.Item(303, 1).Value = 10
$ws.Cells.Item(303, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(303, 3).Value = "La Araucanía"
$ws.Cells.Item(303, 4).Value = 45212
$ws.Cells.Item(303, 5).Value = 9
$ws.Cells.Item(303, 6).Value = 100112005
$ws.Cells.Item(303, 7).Value = "Puerro"
$ws.Cells.Item(303, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(303, 9).Value = "Primera"
$ws.Cells.Item(303, 10).Value = 40
$ws.Cells.Item(303, 11).Value = 8000
$ws.Cells.Item(303, 12).Value = 8000
$ws.Cells.Item(303, 13).Value = 8000
$ws.Cells.Item(303, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(303, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(303, 16).Value = 667
$ws.Cells.Item(303, 17).Value = 12
$ws.Cells.Item(303, 18).Value = "Hortaliza"
